# Added code for Excel reading: populate the "Status" column (E) for the
# second data row with a new "pass" value, mirroring how a test run would
# record its result, and leave the sheet positioned/zoomed on that cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Write the new status value into E2 (creates the new shared string "pass").
$ws.Range("E2").Value = "pass"

# Bring the sheet into view, zoomed in, with the new cell selected -
# matches the updated sheetView/selection state after the edit.
$ws.Activate()
$excel.ActiveWindow.Zoom = 160
$ws.Range("E2").Select()
